$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# This sheet ("Export") is a flat Conta / Nome / Saldo export. The edit:
#
#  1. Deletes the row for account 004550415 (DIOGO, Saldo -251.62)
#  2. Deletes the row for account 005061124 (BRUNO, Saldo -880.6)
#  3. Deletes the existing row for account 004482090 (CEZAR, Saldo 69.96)
#  4. Inserts a brand-new row right under the first data row (new row 3,
#     pushing everything from the old row 3 down by one) for account
#     004482090 (CEZAR) with an updated Saldo of 119069.96
#
# Row numbers are located dynamically with Find() instead of being
# hard-coded, then the two lower rows are removed first (bottom-up) so
# that row numbers found earlier remain valid, and the new row is
# inserted last.
# -----------------------------------------------------------------------

$colA = $ws.Columns.Item(1)

$rowDiogo  = $colA.Find("004550415").Row
$rowBruno  = $colA.Find("005061124").Row
$rowCezarOld = $colA.Find("004482090").Row

$rowsToDelete = @($rowDiogo, $rowBruno, $rowCezarOld) | Sort-Object -Descending
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}

# New row goes immediately after the very first data row (row 2), i.e.
# it becomes the new row 3 and the old row 3 (and below) shift down one.
$ws.Rows.Item(3).Insert()

$newRow = $ws.Rows.Item(3)
$newRow.Cells.Item(1,1).NumberFormat = "@"
$newRow.Cells.Item(1,1).Value2 = "004482090"
$newRow.Cells.Item(1,2).Value2 = "CEZAR"
$newRow.Cells.Item(1,3).Value2 = 119069.96

Write-Host "edit complete"
